$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Update the SUMMARY paragraph: "...will graduate in 2020." ->
#    "...will graduate in the Fall of 2019."
# ---------------------------------------------------------------------------
$rGrad = $d.Content
$okGrad = $rGrad.Find.Execute("will graduate in 2020.", $true, $false, $false, $false, $false, $true, 1, $false, "will graduate in the Fall of 2019.", 2)

# ---------------------------------------------------------------------------
# 2. Update the Texas State University education dates: "(2016 to 2020)" ->
#    "(2016 to 2019)"
# ---------------------------------------------------------------------------
$rDates = $d.Content
$okDates = $rDates.Find.Execute("(2016 to 2020)", $true, $false, $false, $false, $false, $true, 1, $false, "(2016 to 2019)", 2)

# ---------------------------------------------------------------------------
# 3. Move the "_GoBack" bookmark from the end of the document to right after
#    the "Texas State University" heading (tracking the location that was
#    last edited, same as Word does automatically while editing).
# ---------------------------------------------------------------------------
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "Texas State University") {
        $targetPara = $p
    }
}

if ($targetPara -ne $null) {
    $insPoint = $targetPara.Range
    [void]$insPoint.MoveEnd(1, -1)
    $insPoint.Collapse(0)

    # Inserting a bookmark directly at a range collapsed onto a paragraph's
    # trailing boundary is unreliable, so nudge a placeholder character in,
    # anchor the bookmark to it, then remove the placeholder again -- leaving
    # a zero-width bookmark exactly after the heading text.
    $insPoint.InsertAfter("~")
    $markRange = $d.Range($insPoint.Start, $insPoint.Start + 1)
    $d.Bookmarks.Add("_GoBack", $markRange)
    $delRange = $d.Range($insPoint.Start, $insPoint.Start + 1)
    [void]$delRange.Delete()
}

Write-Output "gradReplace=$okGrad datesReplace=$okDates bookmarkMoved=$($targetPara -ne $null)"
